# hit_miss_rule.xlsx -- collapse the two-row header (row 3 = column titles,
# row 4 = "(Type I)"/"(Type II)" sub-titles) into a single header row 4 whose
# labels summarize what used to be spread across both rows, then relabel the
# "Type I + Type II" / "Choice" rows to their new names.
#
# Rows 5-11 (the data rows) keep their row numbers and formulas untouched;
# only row 3 disappears and row 4 becomes the new, single header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rewrite row 4 with the consolidated header text.
#    (Old F3 "Type I + Type II" becomes "Overall Error Rate"; old H3/I3
#    "TI"/"T2" move down one row to H4/I4; the old row-4 sub headers
#    "(Type I)"/"(Type II)" are discarded entirely.)
# ---------------------------------------------------------------------
$ws.Range("C4").Value = "Rule"
$ws.Range("D4").Value = "% incorrectly assigned to control "
$ws.Range("E4").Value = "% incorrectly assigned to treatment"
$ws.Range("F4").Value = "Overall Error Rate"
$ws.Range("H4").Value = "TI"
$ws.Range("I4").Value = "T2"

# ---------------------------------------------------------------------
# 2) Format the new header row: C4:F4 centered with a thin top rule and a
#    double bottom rule (closing out the header block), H4:I4 simply
#    centered like the old H3:I3 cells were.
# ---------------------------------------------------------------------
$hdr = $ws.Range("C4:F4")
$hdr.HorizontalAlignment = -4108
$hdr.Borders.Item(8).LineStyle = 1
$hdr.Borders.Item(9).LineStyle = -4119

$ws.Range("H4:I4").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 3) Row 3 (the old top header row) is no longer needed now that its
#    content lives in row 4 -- clear it in place (no row shifting, so
#    rows 5-11 keep their addresses and formulas).
# ---------------------------------------------------------------------
$ws.Rows.Item(3).Clear()

# ---------------------------------------------------------------------
# 4) Rename the last rule from "Choice" to "Allow choice".
# ---------------------------------------------------------------------
$ws.Range("C10").Value = "Allow choice"

# ---------------------------------------------------------------------
# 5) Update the selection to match the new header-inclusive range.
# ---------------------------------------------------------------------
$ws.Range("C4:F10").Select()
